$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "TO_ASK" header, same formatting as the A1 "TO_DO" header cell
[void]$ws.Range("A1").Copy()
[void]$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "TO_ASK"

# Row 8: the question text, same formatting as a regular body cell (e.g. A2)
[void]$ws.Range("A2").Copy()
[void]$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Czy pasowałoby utworzyć partycje na tabelach głównych oraz relacyjnych na poszczególne okresy czasowe? Czasy zapytań i operacji powinny ulec poprawie"

# Row heights: row 6 is a blank spacer row, row 7 matches the header rows, row 8 is taller to fit wrapped text
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75
$ws.Rows.Item(8).RowHeight = 60.75

$excel.CutCopyMode = $false

# Update the active selection to match the author's last edit position
[void]$ws.Range("B8").Select()
